$d = $word.ActiveDocument

# Replace "<id>p056r_a1</id>" with "<id>p056r_1</id>"
$d.Content.Find.Execute("<id>p056r_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p056r_1</id>", 2)

# Replace "<id>p056r_a2</id>" with "<id>p056r_2</id>"
$d.Content.Find.Execute("<id>p056r_a2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p056r_2</id>", 2)
